$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unmerge existing merged ranges that change ---
$ws.Range("A4:A8").UnMerge()
$ws.Range("A12:A14").UnMerge()
$ws.Range("A17:A18").UnMerge()
$ws.Range("A20:A21").UnMerge()
$ws.Range("A22:A23").UnMerge()

# --- Clear cells that become blank ---
$ws.Range("A9").Value = ""
$ws.Range("A10").Value = ""
$ws.Range("A11").Value = ""
$ws.Range("A12").Value = ""
$ws.Range("A19").Value = ""

# --- Set new / updated cell values ---
$ws.Range("C4").Value = -0.04660037860037861
$ws.Range("E4").Value = 0.04814836814836815
$ws.Range("F4").Value = -0.03184008784008784
$ws.Range("H4").Value = 0.01999341199341199
$ws.Range("I4").Value = 0.02927555327555328
$ws.Range("J4").Value = 0.005216214340192776
$ws.Range("C5").Value = -0.01744488544488544
$ws.Range("E5").Value = -0.03236231636231637
$ws.Range("F5").Value = 0.00459020859020859
$ws.Range("H5").Value = 0.0001692601692601693
$ws.Range("I5").Value = -0.005886209886209887
$ws.Range("J5").Value = -0.0116632175894564
$ws.Range("C6").Value = 0.05373526173526175
$ws.Range("E6").Value = 0.02668350268350269
$ws.Range("F6").Value = -0.02636527436527437
$ws.Range("H6").Value = 0.02171052971052971
$ws.Range("I6").Value = 0.02856575256575257
$ws.Range("J6").Value = -0.01643892812163158
$ws.Range("C7").Value = -0.04066158466158466
$ws.Range("E7").Value = 0.04567797367797368
$ws.Range("F7").Value = 0.04757169557169558
$ws.Range("H7").Value = -0.06304071904071905
$ws.Range("I7").Value = -0.05188736788736789
$ws.Range("J7").Value = -0.04187123295761444
$ws.Range("C8").Value = 0.1192478392478393
$ws.Range("E8").Value = -0.01432170232170232
$ws.Range("F8").Value = -0.9620984420984422
$ws.Range("H8").Value = 0.9999999879999882
$ws.Range("I8").Value = 0.9744730584730585
$ws.Range("J8").Value = -0.02660993208069475
$ws.Range("B9").Value = "Cane  PL [% lipid]"
$ws.Range("C9").Value = 0.009958257958257959
$ws.Range("E9").Value = -0.009452217452217452
$ws.Range("F9").Value = -0.02722975522975523
$ws.Range("H9").Value = 0.03406724206724207
$ws.Range("I9").Value = 0.03017574617574618
$ws.Range("J9").Value = 0.06075327373101429
$ws.Range("B10").Value = "Sorghum  PL [% lipid]"
$ws.Range("C10").Value = -0.001412029412029412
$ws.Range("E10").Value = -0.07043482643482644
$ws.Range("F10").Value = 0.0945033945033945
$ws.Range("H10").Value = -0.07774355374355375
$ws.Range("I10").Value = -0.09023705423705423
$ws.Range("J10").Value = -0.002918144819671102
$ws.Range("B11").Value = "Cane  FFA [% lipid]"
$ws.Range("C11").Value = -0.04223759423759424
$ws.Range("E11").Value = -0.02775066375066376
$ws.Range("F11").Value = 0.01621798021798022
$ws.Range("H11").Value = -0.0082994962994963
$ws.Range("I11").Value = -0.01304689304689305
$ws.Range("J11").Value = -0.0347014160249794
$ws.Range("B12").Value = "Sorghum  FFA [% lipid]"
$ws.Range("C12").Value = -0.03681943281943282
$ws.Range("E12").Value = 0.05713671313671313
$ws.Range("F12").Value = 0.03387094587094588
$ws.Range("H12").Value = -0.05202722802722803
$ws.Range("I12").Value = -0.03650521250521251
$ws.Range("J12").Value = 0.01015733673348077
$ws.Range("B13").Value = "TAG to  FFA conversion [% lipid]"
$ws.Range("C13").Value = 0.0161988881988882
$ws.Range("E13").Value = 0.01318005718005718
$ws.Range("F13").Value = -0.01965318765318766
$ws.Range("H13").Value = 0.01691886491886492
$ws.Range("I13").Value = 0.01803684603684604
$ws.Range("J13").Value = -0.002406602664272236
$ws.Range("A14").Value = "Stream-ethanol"
$ws.Range("B14").Value = "Price [USD/gal]"
$ws.Range("C14").Value = 0.9701300141300142
$ws.Range("E14").Value = -0.03344753744753745
$ws.Range("F14").Value = -0.001511917511917512
$ws.Range("H14").Value = 0.008536532536532537
$ws.Range("I14").Value = 0.002276042276042276
$ws.Range("J14").Value = 0.03768333322412549
$ws.Range("A15").Value = "Stream-biodiesel"
$ws.Range("B15").Value = "Price [USD/gal]"
$ws.Range("C15").Value = -0.01582796782796783
$ws.Range("E15").Value = -0.03454244254244254
$ws.Range("F15").Value = 0.06554944154944156
$ws.Range("H15").Value = -0.05822072222072223
$ws.Range("I15").Value = -0.06354027954027953
$ws.Range("J15").Value = 0.0104322665472174
$ws.Range("A16").Value = "Stream-natural gas"
$ws.Range("B16").Value = "Price [USD/cf]"
$ws.Range("C16").Value = -0.04641826641826642
$ws.Range("E16").Value = -0.04329402729402729
$ws.Range("F16").Value = 0.01197364797364798
$ws.Range("H16").Value = 0.0008034608034608035
$ws.Range("I16").Value = -0.00907874107874108
$ws.Range("J16").Value = 0.07670941131240012
$ws.Range("A17").Value = "biorefinery"
$ws.Range("B17").Value = "Electricity price [USD/kWh]"
$ws.Range("C17").Value = 0.03646461646461646
$ws.Range("E17").Value = -0.01070827070827071
$ws.Range("F17").Value = -0.04543846543846543
$ws.Range("H17").Value = 0.001391053391053391
$ws.Range("I17").Value = 0.001922293922293922
$ws.Range("J17").Value = 0.04806647342139139
$ws.Range("B18").Value = "Operating days [day/yr]"
$ws.Range("C18").Value = 0.02955667755667757
$ws.Range("E18").Value = 0.9999962079962079
$ws.Range("F18").Value = -0.2485104565104566
$ws.Range("H18").Value = -0.01605245205245205
$ws.Range("I18").Value = 0.2044309804309805
$ws.Range("J18").Value = -0.003509900858782483
$ws.Range("B19").Value = "IRR [%]"
$ws.Range("C19").Value = -0.1362950082950083
$ws.Range("E19").Value = -0.02503814503814504
$ws.Range("F19").Value = 0.009814413814413815
$ws.Range("H19").Value = -0.005942513942513943
$ws.Range("I19").Value = -0.01207676407676408
$ws.Range("J19").Value = -0.05990334764380797
$ws.Range("A20").Value = "Stream-crude glycerol"
$ws.Range("B20").Value = "Price [USD/kg]"
$ws.Range("C20").Value = 0.03772921372921374
$ws.Range("E20").Value = -0.02532253332253332
$ws.Range("F20").Value = -0.007780099780099781
$ws.Range("H20").Value = 0.01528855528855529
$ws.Range("I20").Value = 0.01021394221394222
$ws.Range("J20").Value = 0.00237202626492977
$ws.Range("A21").Value = "Stream-pure glycerine"
$ws.Range("B21").Value = "Price [USD/kg]"
$ws.Range("C21").Value = -0.01622202422202422
$ws.Range("E21").Value = 0.05032194232194233
$ws.Range("F21").Value = -0.00677890277890278
$ws.Range("H21").Value = -0.006392730392730393
$ws.Range("I21").Value = 0.004638220638220639
$ws.Range("J21").Value = 0.03955213416468962
$ws.Range("A22").Value = "Stream-cellulase"
$ws.Range("B22").Value = "Price [USD/kg]"
$ws.Range("C22").Value = -0.04204276204276204
$ws.Range("E22").Value = 0.01369310569310569
$ws.Range("F22").Value = -0.02215040215040215
$ws.Range("H22").Value = 0.0211995691995692
$ws.Range("I22").Value = 0.02262506262506263
$ws.Range("J22").Value = -0.01606435246373224
$ws.Range("B23").Value = "Cellulase loading [wt. % cellulose]"
$ws.Range("C23").Value = -0.00171994971994972
$ws.Range("E23").Value = 0.02661846261846262
$ws.Range("F23").Value = -0.02673075873075873
$ws.Range("H23").Value = 0.01831840231840232
$ws.Range("J23").Value = 0.002057788165313365
$ws.Range("A24").Value = "Pretreatment reactor system"
$ws.Range("B24").Value = "Base cost [million USD]"
$ws.Range("C24").Value = -0.05941297141297141
$ws.Range("E24").Value = -0.01913898713898714
$ws.Range("F24").Value = 0.05071770271770272
$ws.Range("H24").Value = -0.04826304026304026
$ws.Range("I24").Value = -0.05215466815466815
$ws.Range("J24").Value = 0.02289902212299082
$ws.Range("A25").Value = "Pretreatment and saccharification"
$ws.Range("B25").Value = "Cane glucose yield [%]"
$ws.Range("C25").Value = -0.01119484719484719
$ws.Range("E25").Value = -0.01535283935283936
$ws.Range("F25").Value = 0.02028609228609229
$ws.Range("H25").Value = -0.01336491736491737
$ws.Range("I25").Value = -0.01586763986763987
$ws.Range("J25").Value = -0.02516321181629988
$ws.Range("B26").Value = "Sorghum glucose yield [%]"
$ws.Range("C26").Value = 0.07816675816675818
$ws.Range("E26").Value = -0.008667824667824668
$ws.Range("F26").Value = 0.0006836886836886837
$ws.Range("H26").Value = 0.002503634503634504
$ws.Range("I26").Value = 0.001016713016713017
$ws.Range("J26").Value = -0.05209387933320162
$ws.Range("B27").Value = "Cane xylose yield [%]"
$ws.Range("C27").Value = -0.001273981273981274
$ws.Range("E27").Value = -0.03613563613563613
$ws.Range("F27").Value = 0.02002238002238002
$ws.Range("H27").Value = -0.01177103977103977
$ws.Range("I27").Value = -0.01690915690915691
$ws.Range("J27").Value = -0.01758809198793372
$ws.Range("B28").Value = "Sorghum xylose yield [%]"
$ws.Range("C28").Value = 0.004481872481872482
$ws.Range("E28").Value = 0.01147054747054747
$ws.Range("F28").Value = 0.02481290481290482
$ws.Range("H28").Value = -0.02956945756945757
$ws.Range("I28").Value = -0.02806369606369606
$ws.Range("J28").Value = -0.02872809383388397
$ws.Range("A29").Value = "Cofermenation"
$ws.Range("B29").Value = "Glucose to ethanol yield [%]"
$ws.Range("C29").Value = 0.007437067437067437
$ws.Range("E29").Value = -0.03008037008037008
$ws.Range("F29").Value = 0.005602121602121602
$ws.Range("H29").Value = -0.0005643005643005644
$ws.Range("I29").Value = -0.007856875856875858
$ws.Range("J29").Value = 0.0310090192731061
$ws.Range("B30").Value = "Xylose to ethanol yield [%]"
$ws.Range("C30").Value = -0.006611106611106611
$ws.Range("E30").Value = 0.0076985236985237
$ws.Range("F30").Value = -0.02701040701040701
$ws.Range("H30").Value = 0.02384813984813985
$ws.Range("I30").Value = 0.02671746271746272
$ws.Range("J30").Value = -0.04849618036139617

# --- Re-merge updated ranges ---
$ws.Range("A4:A13").Merge()
$ws.Range("A17:A19").Merge()
$ws.Range("A22:A23").Merge()
$ws.Range("A25:A28").Merge()
$ws.Range("A29:A30").Merge()